$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 2.3
$ws.Range("H2").Value = 2.65
$ws.Range("I2").Value = 3.7
$ws.Range("J2").Value = 3.25
$ws.Range("M2").Value = 1.17
$ws.Range("N2").Value = 5
$ws.Range("X2").Value = 9.5
$ws.Range("AW2").Value = 5
# Row 3
$ws.Range("G3").Value = 2.7
$ws.Range("H3").Value = 2.75
$ws.Range("I3").Value = 2.85
# Row 4
$ws.Range("M4").Value = 1.05
$ws.Range("N4").Value = 6.5
$ws.Range("O4").Value = 1.47
$ws.Range("P4").Value = 2.37
# Row 5
$ws.Range("M5").Value = 1.02
$ws.Range("O5").Value = 1.22
# Row 6
$ws.Range("G6").Value = 2.63
$ws.Range("I6").Value = 2.63
$ws.Range("J6").Value = 3.2
$ws.Range("M6").Value = 1.02
$ws.Range("O6").Value = 1.19
$ws.Range("S6").Value = 1.36
$ws.Range("T6").Value = 3
$ws.Range("U6").Value = 1.62
$ws.Range("V6").Value = 2.2
$ws.Range("W6").Value = 10
$ws.Range("AB6").Value = 26
$ws.Range("AC6").Value = 12
$ws.Range("AE6").Value = 12
$ws.Range("AG6").Value = 10
$ws.Range("AL6").Value = 26
$ws.Range("AM6").Value = 151
$ws.Range("AP6").Value = 21
$ws.Range("AQ6").Value = 41
$ws.Range("AR6").Value = 51
$ws.Range("AT6").Value = 3
$ws.Range("AY6").Value = 21
$ws.Range("BA6").Value = 51
# Row 10
$ws.Range("G10").Value = 2.45
$ws.Range("H10").Value = 2.65
# Row 12
$ws.Range("G12").Value = 1.29
# Row 13
$ws.Range("K13").Value = 2.05
$ws.Range("L13").Value = 5
$ws.Range("S13").Value = 1.42
$ws.Range("T13").Value = 2.45
$ws.Range("W13").Value = 5.7
$ws.Range("X13").Value = 7.4
$ws.Range("Y13").Value = 8.25
$ws.Range("Z13").Value = 14
$ws.Range("AB13").Value = 32
$ws.Range("AF13").Value = 100
$ws.Range("AO13").Value = 8.5
$ws.Range("AP13").Value = 19
$ws.Range("AQ13").Value = 30
$ws.Range("AR13").Value = 65
$ws.Range("AS13").Value = 250
$ws.Range("AT13").Value = 2.4
$ws.Range("AU13").Value = 7.6
